$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" (column C) values for rows 2-252 in three groups,
# matching the change described in the diff.

# Rows 2-8: 7310 -> 7623
$ws.Range("C2:C8").Value = 7623

# Rows 9-46: 7310 -> 7590
$ws.Range("C9:C46").Value = 7590

# Rows 47-252: (7310 / 7295 / 7293) -> 7573
$ws.Range("C47:C252").Value = 7573
